$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert two new rows at row 13 (shifts old rows 13-23 down to 15-25)
$ws.Rows.Item(13).Resize(2).Insert()

# 2) Apply formatting (style) to the newly created B13/C13 and B14/C14 cells
#    by copying formats from the existing B/C cells in row 15 (which still has
#    the original style-2/style-3 formatting used throughout column B/C).
$ws.Range("B15").Copy()
$ws.Range("B13:B14").PasteSpecial(-4122)
$ws.Range("C15").Copy()
$ws.Range("C13:C14").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 3) Set the new/updated cell values
$ws.Range("B10").Value = 'Passar aos alunos os conhecimentos básicos da estrutura de um vidro, a influência da composição nas propriedades e os processos empregados na produção de vidros'
$ws.Range("C10").Value = 'Passar aos alunos os conhecimentos básicos da estrutura de um vidro, a influência da composição nas propriedades e os processos empregados na produção de vidros'
$ws.Range("B13").Value = '5983729 - Fernando Vernilli Junior'
$ws.Range("C13").Value = '5983729 - Fernando Vernilli Junior'
$ws.Range("B14").Value = '1922320 - Sebastiao Ribeiro'
$ws.Range("C14").Value = '1922320 - Sebastiao Ribeiro'
$ws.Range("B15").Value = '1.Introdução, Quadro da indústria brasileira de vidros2.Composição dos vidros3.Materiais Primas4.Mecanismo de fusão e formação do vidro5.Viscosidade - Definição, relação com a composição, métodos experimentais de medição, cálculo a partir da composição6.Propriedades óticas 7.Propriedades mecânicas 8.Propriedades químicas 9.Processamento - Vidro plano, vidro oco, vidros especiais, vidro temperado, esmalte10.Aula prática - Fundir um vidro, produzir um vidro colorido e esmaltar um metal'
$ws.Range("C15").Value = '1.Introdução, Quadro da indústria brasileira de vidros2.Composição dos vidros3.Materiais Primas4.Mecanismo de fusão e formação do vidro5.Viscosidade - Definição, relação com a composição, métodos experimentais de medição, cálculo a partir da composição6.Propriedades óticas 7.Propriedades mecânicas 8.Propriedades químicas 9.Processamento - Vidro plano, vidro oco, vidros especiais, vidro temperado, esmalte10.Aula prática - Fundir um vidro, produzir um vidro colorido e esmaltar um metal'
$ws.Range("B17").Value = '1. Introdução, Quadro da indústria brasileira de vidros2. Composição dos vidros3. Materiais Primas4. Mecanismo de fusão e formação do vidro5. Viscosidade – Definição, relação com a composição, métodos experimentais de medição, cálculo a partir da composição6. Propriedades óticas 7. Propriedades mecânicas 8. Propriedades químicas 9. Processamento – Vidro plano, vidro oco, vidros especiais, vidro temperado, esmalte10. Aula prática - Fundir um vidro, produzir um vidro colorido e esmaltar um metal'
$ws.Range("C17").Value = '1. Introdução, Quadro da indústria brasileira de vidros2. Composição dos vidros3. Materiais Primas4. Mecanismo de fusão e formação do vidro5. Viscosidade – Definição, relação com a composição, métodos experimentais de medição, cálculo a partir da composição6. Propriedades óticas 7. Propriedades mecânicas 8. Propriedades químicas 9. Processamento – Vidro plano, vidro oco, vidros especiais, vidro temperado, esmalte10. Aula prática - Fundir um vidro, produzir um vidro colorido e esmaltar um metal'
$ws.Range("B20").Value = 'Serão realizadas duas provas escritas (P1 e P2), apresentações orais de trabalhos (T) e listas de exercícios (E)'
$ws.Range("C20").Value = 'Serão realizadas duas provas escritas (P1 e P2), apresentações orais de trabalhos (T) e listas de exercícios (E)'
$ws.Range("B21").Value = 'A nota final será calculada utilizando a equação: {[(P1 + P2 + T)/3] x 0,9} + E x 0,1'
$ws.Range("C21").Value = 'A nota final será calculada utilizando a equação: {[(P1 + P2 + T)/3] x 0,9} + E x 0,1'
$ws.Range("B22").Value = 'Para a recuperação será realizada uma prova  (PR) abrangendo toda a matéria no semestre, valendo de 0 (zero) a 10 (10). Média Final: (MP + PR)/2. Média Final igual ou superior a 5: aprovado. Média Final inferior a 5: reprovado'
$ws.Range("C22").Value = 'Para a recuperação será realizada uma prova  (PR) abrangendo toda a matéria no semestre, valendo de 0 (zero) a 10 (10). Média Final: (MP + PR)/2. Média Final igual ou superior a 5: aprovado. Média Final inferior a 5: reprovado'
$ws.Range("B23").Value = '1.)Associação Brasileira da Industria de Vidros, www.abividro.br2.)H. Scholze, Glas, Springer-Verlag, 19883.)R. H. Doremus, Glass Science, New York, John Wiley, 19944.)H. G. Pfaender, Schott Guide to Glass, London, Chapman & Hall, 1996'
$ws.Range("C23").Value = '1.)Associação Brasileira da Industria de Vidros, www.abividro.br2.)H. Scholze, Glas, Springer-Verlag, 19883.)R. H. Doremus, Glass Science, New York, John Wiley, 19944.)H. G. Pfaender, Schott Guide to Glass, London, Chapman & Hall, 1996'

# 4) Clean up the column definitions: column A (1) should only be 30.71 wide
#    on its own; column B (2) keeps its own 60.71 width definition. Nudging
#    column B's width causes the engine to split the old merged A:B column
#    range, isolating column A's width entry.
$ws.Columns.Item(2).ColumnWidth = 60.7109375
